# Update the "dSF" column (F) values to reflect the repulled/pushed data
# and recalculated mean, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -2
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -2
$ws.Range("F10").Value = -4
$ws.Range("F14").Value = -1
$ws.Range("F18").Value = -3
$ws.Range("F23").Value = 5
$ws.Range("F27").Value = -3
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = -9
